$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = "left, pin 2"
$ws.Range("D12").Value = "right, pin 2"

$ws.Range("C11").Value = "QEI0"
$ws.Range("C12").Value = "QEI1"

$ws.Range("A12").Value = "C"
$ws.Range("A11").Value = "F"

$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 5

$ws.Range("B12").Select()
